# Update NATMI LR-pair worksheet (Wnt7b-Fzd4) with newly recomputed TPM values.
# - Recomputed statistics for existing sender/target rows (rows 2-5, FAPs sender)
# - Replaced the "Resolving-Mac" target-cluster rows with recomputed rows for the
#   MuSCs sender (rows 6-9), since the "Resolving-Mac" cluster no longer exists
#   in the new TPM run
# - Removed the now-obsolete trailing rows (old rows 10-11)
# - The unused "Resolving-Mac" shared string is dropped automatically because no
#   cell references it any more once the sheet is rewritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final data for rows 2-9 (columns A:T).
$data = @(
    ,@("FAPs","Wnt7b","Fzd4","ECs", 1, 0.3333333333333333, 0.021087, 0.063261, 0.02328126719340038, 0.02328126719340038, 3, 1, 26.60444266666667, 79.813328, 0.5736225649467147, 0.5736225649467147, 0.561007882512, 5.049070942608, 0.01335466020268813, 0.01335466020268813)
    ,@("FAPs","Wnt7b","Fzd4","FAPs", 1, 0.3333333333333333, 0.021087, 0.063261, 0.02328126719340038, 0.02328126719340038, 3, 1, 14.03147833333333, 42.094435, 0.3025349071358453, 0.3025349071358453, 0.295881783615, 2.662936052535, 0.007043396008360186, 0.007043396008360187)
    ,@("FAPs","Wnt7b","Fzd4","Inflammatory-Mac", 1, 0.3333333333333333, 0.021087, 0.063261, 0.02328126719340038, 0.02328126719340038, 1, 0.3333333333333333, 0.2022123333333333, 0.606637, 0.004359931864156574, 0.004359931864156574, 0.004264051473, 0.038376463257, 0.0001015047386744494, 0.0001015047386744494)
    ,@("FAPs","Wnt7b","Fzd4","MuSCs", 1, 0.3333333333333333, 0.021087, 0.063261, 0.02328126719340038, 0.02328126719340038, 3, 1, 5.541567000000001, 16.624701, 0.1194825960532834, 0.1194825960532834, 0.116855023329, 1.051695209961, 0.002781706243677616, 0.002781706243677617)
    ,@("MuSCs","Wnt7b","Fzd4","ECs", 3, 1, 0.8846626666666667, 2.653988, 0.9767187328065996, 0.9767187328065997, 3, 1, 26.60444266666667, 79.813328, 0.5736225649467147, 0.5736225649467147, 23.53595719467378, 211.823614752064, 0.5602679047440265, 0.5602679047440267)
    ,@("MuSCs","Wnt7b","Fzd4","FAPs", 3, 1, 0.8846626666666667, 2.653988, 0.9767187328065996, 0.9767187328065997, 3, 1, 14.03147833333333, 42.094435, 0.3025349071358453, 0.3025349071358453, 12.41312503964222, 111.71812535678, 0.2954915111274851, 0.2954915111274851)
    ,@("MuSCs","Wnt7b","Fzd4","Inflammatory-Mac", 3, 1, 0.8846626666666667, 2.653988, 0.9767187328065996, 0.9767187328065997, 1, 0.3333333333333333, 0.2022123333333333, 0.606637, 0.004359931864156574, 0.004359931864156574, 0.1788897020395556, 1.610007318356, 0.004258427125482125, 0.004258427125482125)
    ,@("MuSCs","Wnt7b","Fzd4","MuSCs", 3, 1, 0.8846626666666667, 2.653988, 0.9767187328065996, 0.9767187328065997, 3, 1, 5.541567000000001, 16.624701, 0.1194825960532834, 0.1194825960532834, 4.902417439732001, 44.121756957588, 0.1167008898096058, 0.1167008898096058)
)

# Remove the two trailing rows (old rows 10 and 11) so the table ends at row 9.
$ws.Rows.Item(11).Delete() | Out-Null
$ws.Rows.Item(10).Delete() | Out-Null

# Write the recomputed values into rows 2-9, columns A-T, one row at a time
# using a 2-D COM array (required for bulk Range.Value assignment).
for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $rowVals = $data[$i]
    $arr = New-Object 'object[,]' 1,20
    for ($j = 0; $j -lt 20; $j++) {
        $arr[0, $j] = $rowVals[$j]
    }
    $ws.Range($ws.Cells.Item($r, 1), $ws.Cells.Item($r, 20)).Value = $arr
}
